$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'57.772.42"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3
$ws.Range("D3").Value = "'3.142.41"
$ws.Range("E3").Value = "  +1.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "'531.27"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("D6").Value = "'140.38"
$ws.Range("E6").Value = "  -0.86%  "

# Row 7
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("D8").Value = "'3.143.20"
$ws.Range("E8").Value = "  +1.11%  "

# Row 9
$ws.Range("D9").Value = "'0.447"
$ws.Range("E9").Value = "  +2.74%  "

# Row 10
$ws.Range("D10").Value = "'7.21"
$ws.Range("E10").Value = "  -1.18%  "

# Row 11
$ws.Range("D11").Value = "'0.109"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("D12").Value = "'0.398"
$ws.Range("E12").Value = "  +3.97%  "

# Row 13
$ws.Range("D13").Value = "'3.693.74"
$ws.Range("E13").Value = "  +1.43%  "

# Row 14
$ws.Range("D14").Value = "'0.134"
$ws.Range("E14").Value = "  +2.94%  "

# Row 15
$ws.Range("D15").Value = "'25.57"
$ws.Range("E15").Value = "  -2.41%  "

# Row 16
$ws.Range("D16").Value = "'0.0000165"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.152.99"
$ws.Range("E17").Value = "  +1.50%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'57.948.57"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19
$ws.Range("D19").Value = "'6.12"
$ws.Range("E19").Value = "  +0.16%  "

# Row 20
$ws.Range("D20").Value = "'12.80"
$ws.Range("E20").Value = "  -0.22%  "

# Row 21
$ws.Range("D21").Value = "'7.97"
$ws.Range("E21").Value = "  -1.09%  "

# Row 22
$ws.Range("D22").Value = "'353.27"
$ws.Range("E22").Value = "  +4.69%  "

# Row 23
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.37%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'68.43"
$ws.Range("E24").Value = "  +3.14%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.509"
$ws.Range("E25").Value = "  +0.46%  "

# Row 26
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = "  +1.26%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0938"
$ws.Range("E28").Value = "  +0.75%  "

# Row 29
$ws.Range("D29").Value = "'7.45"
$ws.Range("E29").Value = "  +3.18%  "

# Row 30
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "'6.39"
$ws.Range("E31").Value = "  -2.82%  "

# Row 32
$ws.Range("D32").Value = "'1.90"
$ws.Range("E32").Value = "  +1.81%  "

# Row 33
$ws.Range("D33").Value = "'21.20"
$ws.Range("E33").Value = "  +1.45%  "

# Row 34
$ws.Range("D34").Value = "'1.21"
$ws.Range("E34").Value = "  +0.34%  "

# Row 35
$ws.Range("D35").Value = "'4.86"
$ws.Range("E35").Value = "  +5.59%  "

# Row 36
$ws.Range("D36").Value = "'157.54"
$ws.Range("E36").Value = "  +2.03%  "

# Row 37
$ws.Range("D37").Value = "'6.15"
$ws.Range("E37").Value = "  +1.54%  "

# Row 38
$ws.Range("D38").Value = "'26.09"
$ws.Range("E38").Value = "  -3.32%  "

# Row 39
$ws.Range("D39").Value = "'1.29"
$ws.Range("E39").Value = "  -0.76%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0671"
$ws.Range("E40").Value = "  +0.73%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.62"
$ws.Range("E41").Value = "  +10.92%  "

# Row 42
$ws.Range("D42").Value = "'4.06"
$ws.Range("E42").Value = "  +4.05%  "

# Row 43
$ws.Range("D43").Value = "'0.703"
$ws.Range("E43").Value = "  +2.92%  "

# Row 44
$ws.Range("D44").Value = "'3.189.35"
$ws.Range("E44").Value = "  +1.21%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'36.65"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0272"
$ws.Range("E46").Value = "  +5.50%  "

# Row 47
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +0.17%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'2.337.40"
$ws.Range("E48").Value = "  +2.64%  "

# Row 49
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  +2.95%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.03"
$ws.Range("E50").Value = "  +0.26%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'20.36"
$ws.Range("E51").Value = "  -1.34%  "
